# Add a new "Italy" worksheet, cloned from the existing "Slovakia" sheet
# (same layout/formatting), positioned right after it, with Italy-specific
# values, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")

# Deselect everything on Slovakia before leaving it, matching the
# "whole sheet" selection state left behind on inactive tabs.
$slovakia.Activate() | Out-Null
$slovakia.Cells.Select() | Out-Null

# Clone the Slovakia sheet (keeps styles, merges, column widths, etc.)
# and place the copy immediately after it.
$slovakia.Copy($null, $slovakia)

$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Country-specific content
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2454/T2453"

# Leave the same relative selection/active cell the template sheets use
$italy.Range("B4").Select() | Out-Null

$italy.Activate() | Out-Null
